$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 (2025-T1) figures
$ws.Range("C3").Value = 18676
$ws.Range("D3").Value = 22715
$ws.Range("E3").Value = 17.78120184899846
$ws.Range("F3").Value = 82.218798151001536

# Update row 4 (2025-T2) figures
$ws.Range("B4").Value = 2890
$ws.Range("C4").Value = 13830
$ws.Range("D4").Value = 16720
$ws.Range("E4").Value = 17.284688995215308
$ws.Range("F4").Value = 82.715311004784681

# Update the active selection shown in the saved view (G7 single cell)
$ws.Range("G7").Select()
